# Bug fix on TSP. log save.
# Applies the data + header changes described by the commit diff to the
# single worksheet in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header block (rows 1-3): two algorithm-label columns were renamed from
# "ACBR_iBNO"/"ACBR_iFAO" + "TSR_LS"/"TSR"/"TSR_RM"/"TSR_RMIE" pairs to a
# smaller, five-column "ACBR_BNO" / "TSR_TWW" / "BTS_TWW" block; the old
# K/L/M columns (TSR_RM / TSR_RMIE pair) no longer carry a sub-label.
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "ACBR_BNO"
$ws.Range("F3").Value = "TSR_TWW"

$ws.Range("G2").Value = "ACBR_BNO"
$ws.Range("G3").Value = "BTS_TWW"

$ws.Range("H2").Value = "ACBR_BNO"
$ws.Range("H3").Value = "TSR_TWW"

$ws.Range("I2").Value = "ACBR_BNO"
# I3 stays "TSR" - unchanged

$ws.Range("J2").Value = "ACBR_BNO"
$ws.Range("J3").Value = "TSR_TWW"

# K2/K3, L2/L3, M2/M3 no longer have sub-headers - clear their contents
# while leaving the (merged/bordered) cell formatting untouched.
$ws.Range("K2").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("M3").ClearContents()

# ---------------------------------------------------------------------
# Data block: updated raw results in columns I and J (rows 6-19, 21-34).
# All other cells on the sheet (N/O MIN columns and the row4/5/20 AVERAGE
# summary rows) are formulas and recompute automatically on recalc.
# ---------------------------------------------------------------------
$ws.Range("I6").Value = 5295
$ws.Range("J6").Value = 4810
$ws.Range("I7").Value = 4430
$ws.Range("J7").Value = 5065
$ws.Range("I8").Value = 3375
$ws.Range("J8").Value = 3175
$ws.Range("I9").Value = 5790
$ws.Range("J9").Value = 5965
$ws.Range("I10").Value = 2135
$ws.Range("J10").Value = 2285
$ws.Range("I11").Value = 3860
$ws.Range("J11").Value = 4065
$ws.Range("J12").Value = 4210
$ws.Range("I13").Value = 6915
$ws.Range("J13").Value = 6955
$ws.Range("I14").Value = 3155
$ws.Range("J14").Value = 3205
$ws.Range("I15").Value = 4155
$ws.Range("J15").Value = 4180
$ws.Range("I16").Value = 1800
$ws.Range("J16").Value = 1840
$ws.Range("I17").Value = 3105
$ws.Range("J17").Value = 3050
$ws.Range("I18").Value = 2040
$ws.Range("J18").Value = 2000
$ws.Range("I19").Value = 2125
$ws.Range("J19").Value = 2035
$ws.Range("I21").Value = 5505
$ws.Range("J21").Value = 5665
$ws.Range("I22").Value = 4700
$ws.Range("J22").Value = 5175
$ws.Range("I23").Value = 3230
$ws.Range("J23").Value = 3270
$ws.Range("I24").Value = 5990
$ws.Range("J24").Value = 6600
$ws.Range("I25").Value = 2580
$ws.Range("J25").Value = 2905
$ws.Range("I26").Value = 4455
$ws.Range("J26").Value = 4990
$ws.Range("I27").Value = 4150
$ws.Range("J27").Value = 3970
$ws.Range("I28").Value = 6745
$ws.Range("J28").Value = 6835
$ws.Range("I29").Value = 3380
$ws.Range("J29").Value = 3435
$ws.Range("I30").Value = 3610
$ws.Range("J30").Value = 3755
$ws.Range("I31").Value = 1810
$ws.Range("J31").Value = 1895
$ws.Range("I32").Value = 2680
$ws.Range("J32").Value = 2735
$ws.Range("I33").Value = 2120
$ws.Range("J33").Value = 2185
$ws.Range("I34").Value = 2215
$ws.Range("J34").Value = 2175

# ---------------------------------------------------------------------
# Selection cursor moved to M10 (was I11).
# ---------------------------------------------------------------------
$ws.Range("M10").Select()
